# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(4, 1, 0, 2, 1, 2, 3, 1, 0, 0, 1, 1, 0, 0, 1, 2, 0, 1, 3, 1, 1, 1, 1, 2, 1, 4, 0, 3, 1, 1, 0, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
